# Expand the 3-column transition table into 5 columns by splitting a
# cell in every row into three cells, then writing the final text for
# every cell in the row. This avoids introducing a <w:tblGrid> (which
# Table.Columns.Add would add) and keeps each cell's paragraph using
# the "Normal" style, matching the rest of the table.

$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# Final text for every row, left-to-right, after the edit (5 columns).
$rows = @(
    @("week4",        "4-&gt;8",      "week8",        "8-&gt;12",     "week12"),
    @("52 (58.43%)",  "30 (34.48%)",  "40 (45.45%)",  "20 (23.53%)",  "26 (30.23%)"),
    @("0 (0%)",       "9 (10.34%)",   "0 (0%)",       "6 (7.06%)",    "0 (0%)"),
    @("37 (41.57%)",  "20 (22.99%)",  "48 (54.55%)",  "17 (20%)",     "60 (69.77%)"),
    @("0 (0%)",       "28 (32.18%)",  "0 (0%)",       "42 (49.41%)",  "0 (0%)"),
    @("89 (100%)",    "87 (100%)",    "88 (100%)",    "85 (100%)",    "86 (100%)")
)

# Column (1-based) whose cell gets split into 3 in each row: the
# header row's middle cell ("week8") sits between the two new cells,
# while every data row's trailing cell is split to grow the new
# columns after the existing two.
$splitColumn = @(2, 3, 3, 3, 3, 3)

for ($r = 1; $r -le $t.Rows.Count; $r++) {
    $col = $splitColumn[$r - 1]
    $cell = $t.Cell($r, $col)
    $cell.Split(1, 2) | Out-Null
    $cell2 = $t.Cell($r, $col + 1)
    $cell2.Split(1, 2) | Out-Null

    $values = $rows[$r - 1]
    for ($c = 1; $c -le 5; $c++) {
        $target = $t.Cell($r, $c)
        $target.Range.Text = $values[$c - 1]
        $target.Range.Style = "Normal"
    }
}
